$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 2 (Sending cluster = ECs) entirely, keeping only the
# former row 3 (Sending cluster = FAPs) which becomes the new (only) data row.
$ws.Rows("2").Delete()

# Update the (now) row 2 values to reflect the new TPM-derived numbers.
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Fzd10"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3532066666666667
$ws.Range("H2").Value = 1.05962
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04795666666666667
$ws.Range("N2").Value = 0.14387
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.01693861437777778
$ws.Range("R2").Value = 0.1524475294
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
